$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.318.24"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.873.59"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "0.7116"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "241.71"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.3109"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "0.07772"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("D10").Value = "25.07"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("D11").Value = "0.08399"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "1.875.47"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "5.244"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "0.7118"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "91.13"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "29.323.44"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "6.062"
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "0.000008190"
$ws.Range("E18").Value = "  +5.04%  "
$ws.Range("D19").Value = "239.57"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").Value = "13.21"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").Value = "2.119.78"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "7.766"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "0.1585"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "162.86"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "9.026"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "4.402"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "4.325"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").Value = "1.286"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("D33").Value = "0.05292"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "1.937"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").Value = "0.7446"
$ws.Range("E36").Value = "  -6.60%  "
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "0.01878"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").Value = "1.219.34"
$ws.Range("E39").Value = "  +4.87%  "
$ws.Range("D40").Value = "2.726"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "6.438"
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("D42").Value = "109.77"
$ws.Range("E42").Value = "  +7.16%  "
$ws.Range("D43").Value = "0.8852"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "72.48"
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "2.017.66"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").Value = "9.368"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "0.4309"
$ws.Range("E51").Value = "  +1.10%  "

# Row 47/48 swap: RenderToken/Mantle order swapped with updated prices
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.5197"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.796"
$ws.Range("E48").Value = "  +1.23%  "
